$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet1)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 2952
$ws.Range("J3").Value = 3057
$ws.Range("H4").Value = 1694
$ws.Range("I4").Value = 1757
$ws.Range("J4").Value = 689
$ws.Range("J5").Value = 235
$ws.Range("J6").Value = 3695
$ws.Range("H7").Value = 26004
$ws.Range("I7").Value = 26204
$ws.Range("J7").Value = 10628

# By Neighborhood (sheet2)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J8").Value = 680
$ws.Range("J9").Value = 64
$ws.Range("J10").Value = 64
$ws.Range("J11").Value = 151
$ws.Range("J19").Value = 338
$ws.Range("I20").Value = 639
$ws.Range("J20").Value = 216
$ws.Range("J23").Value = 110
$ws.Range("J26").Value = 16
$ws.Range("J29").Value = 601
$ws.Range("J30").Value = 44
$ws.Range("J33").Value = 446
$ws.Range("J34").Value = 55
$ws.Range("J37").Value = 348
$ws.Range("J42").Value = 427
$ws.Range("J44").Value = 82
$ws.Range("J45").Value = 14
$ws.Range("J49").Value = 68
$ws.Range("J52").Value = 281
$ws.Range("J55").Value = 135
$ws.Range("J60").Value = 68
$ws.Range("J63").Value = 51
$ws.Range("J65").Value = 281
$ws.Range("J67").Value = 380
$ws.Range("J72").Value = 39
$ws.Range("J75").Value = 33
$ws.Range("J76").Value = 153
$ws.Range("J77").Value = 95
$ws.Range("J78").Value = 138
$ws.Range("J79").Value = 314
$ws.Range("J83").Value = 247
$ws.Range("J85").Value = 489
$ws.Range("J89").Value = 116
$ws.Range("J90").Value = 122
$ws.Range("J93").Value = 53
$ws.Range("H95").Value = 513
$ws.Range("J95").Value = 171
$ws.Range("J96").Value = 122
$ws.Range("I98").Value = 186
$ws.Range("J99").Value = 154
$ws.Range("H101").Value = 26004
$ws.Range("I101").Value = 26204
$ws.Range("J101").Value = 10628

# South Shore (sheet3)
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 118
$ws.Range("J3").Value = 186
$ws.Range("J6").Value = 137
$ws.Range("J7").Value = 489

# Little Village (sheet5)
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 64
$ws.Range("J3").Value = 75
$ws.Range("J6").Value = 125
$ws.Range("J7").Value = 281

# Belmont Cragin (sheet6)
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 151

# Austin (sheet7)
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 205
$ws.Range("J3").Value = 218
$ws.Range("J6").Value = 202
$ws.Range("J7").Value = 680

# Uptown (sheet10)
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 116

# West Ridge (sheet11)
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 38
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 122

# Fuller Park (sheet13)
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J2").Value = 18
$ws.Range("J3").Value = 14
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 44

# Grand Crossing (sheet14)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 105
$ws.Range("J7").Value = 348

# Woodlawn (sheet15)
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 154

# North Lawndale (sheet16)
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J4").Value = 25
$ws.Range("J7").Value = 380

# New City (sheet19)
$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 82
$ws.Range("J7").Value = 281

# South Chicago (sheet20)
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 90
$ws.Range("J7").Value = 247

# West Pullman (sheet21)
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("H4").Value = 19
$ws.Range("J4").Value = 6
$ws.Range("H7").Value = 513
$ws.Range("J7").Value = 171

# Garfield Park (sheet22)
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 118
$ws.Range("J6").Value = 151
$ws.Range("J7").Value = 446

# Lincoln Park (sheet23)
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J3").Value = 14
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 68

# Loop (sheet24)
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 40
$ws.Range("J6").Value = 96

# Englewood (sheet25)
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 177
$ws.Range("J7").Value = 601

# Chatham (sheet26)
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 80
$ws.Range("J3").Value = 94
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 338

# Irving Park (sheet27)
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J3").Value = 22
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 82

# River North (sheet29)
$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 153

# Humboldt Park (sheet32)
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 86
$ws.Range("J3").Value = 91
$ws.Range("J6").Value = 219
$ws.Range("J7").Value = 427

# Avondale (sheet34)
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 64

# Rogers Park (sheet35)
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 138

# Lower West Side (sheet36)
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 27
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 135

# Douglas (sheet39)
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 110

# Roseland (sheet42)
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 89
$ws.Range("J3").Value = 114
$ws.Range("J6").Value = 88
$ws.Range("J7").Value = 314

# Chicago Lawn (sheet44)
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 72
$ws.Range("I4").Value = 41
$ws.Range("I7").Value = 639
$ws.Range("J7").Value = 216

# West Lawn (sheet48)
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 53

# Garfield Ridge (sheet50)
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 55

# Wicker Park (sheet55)
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 186

# Lincoln Square (sheet56)
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J3").Value = 19
$ws.Range("J6").Value = 15

# East Village (sheet57)
$ws = $wb.Worksheets.Item('East Village')
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 16

# Avalon Park (sheet61)
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 64

# Pullman (sheet73)
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 33

# Washington Heights (sheet74)
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 41
$ws.Range("J3").Value = 37
$ws.Range("J7").Value = 122

# Morgan Park (sheet78)
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 68

# Old Town (sheet82)
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 39

# Riverdale (sheet84)
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 31
$ws.Range("J7").Value = 95

# Jackson Park (sheet85)
$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 14
